# Generate Report for Handoff
# Updates the "latest handoff / generate" timestamps for the
# f3feb4d7-776b-4f2d-adf2-2d4b76f883fb file (row 5 on every sheet)
# to reflect a freshly generated handoff report.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: "Latest Handoff Datetime" column (H) for row 5
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H5").Value = "2017-02-22 06:42:10"

# de-de sheet: "Latest Handoff Datetime" column (H) for row 5
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H5").Value = "2017-02-22 06:42:26"

# Overview sheet: "Latest HO Xliff Generate Date" column (G) for row 5
# reflects the most recent of the per-language handoff datetimes above.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G5").Value = "2017-02-22 06:42:26"
